$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the bold "Differences" heading paragraph and replace it
#    with the new (non-bold) narrative paragraphs.
# ------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Differences") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Differences' paragraph"
}

$newParagraphs = @(
    "We have talked about related micropayment systems. This thesis is related with an earlier work. Can Yücel, improved these payment systems and adopted them for a real time support. A real time system with high frequency of payments, which basically means a lot of payments per hour.",
    "I took over Can’s work and improved it to support some new features.",
    "These features are anonymity, untraceability, seamless roaming.",
    "Can conducted his simulations using a network simulator called OMNET++. This simulator is very popular and neat but it does not support mesh networks. Can simulated ad hoc networks with IEEE 802.11b/g.",
    "We have implemented our system in ns-3 from scratch.",
    "Can’s simulations had burst scenarios and low demand scenarios; they were not really close to real life situations.",
    "I have implemented client types, which affects the frequency of network usage and mobility speed and distance."
)

# Turn off the bold heading formatting and drop in the first paragraph's text.
$r = $target.Range
$r.Font.Bold = $false
$r.Text = $newParagraphs[0]

# Append the remaining paragraphs right after, inheriting the (now
# non-bold) paragraph formatting of the first one.
$prev = $target
for ($i = 1; $i -lt $newParagraphs.Length; $i++) {
    $prev.Range.InsertParagraphAfter()
    $prev = $prev.Next()
    $prev.Range.Text = $newParagraphs[$i]
}

Write-Host "Replaced 'Differences' heading with $($newParagraphs.Length) narrative paragraphs."
